# card finalizado+bt add/cancel + logica
# Updates row 10 (Claudio Castro Monstro) with computed Idade/Categoria/Data de
# Nascimento, and appends three new "chamada" rows (11-13) coming from the
# add-student card - two fully processed by the age/category logic (11, 12)
# and a third (13) added but not yet run through that logic, mirroring the
# "cancel" partial state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alunos")

# Date cells in this sheet use a YYYY-MM-DD HH:MM:SS number format (style
# already present on column M for rows 2-9); reuse it so we don't fork a new
# style entry when writing the new birth-date serials below.
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 10: finish the pending "Claudio Castro Monstro" entry ------------
$ws.Cells.Item(10, 9).Value = 46                 # I10 Idade
$ws.Cells.Item(10, 10).Value = "Não definida"     # J10 Categoria
$ws.Cells.Item(10, 13).NumberFormat = $dateFormat
$ws.Cells.Item(10, 13).Value = 28897              # M10 Data de Nascimento

# --- Row 11: julio cesar darwin (fully processed) --------------------------
$ws.Cells.Item(11, 1).Value = "julio cesar darwin"
$ws.Cells.Item(11, 2).Value = ""
$ws.Cells.Item(11, 3).Value = "Sim"
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = "Terça e Quinta"
$ws.Cells.Item(11, 6).Value = "16h00"
$ws.Cells.Item(11, 7).Value = "Jefferson"
$ws.Cells.Item(11, 8).Value = "Adulto B"
$ws.Cells.Item(11, 9).Value = 24
$ws.Cells.Item(11, 10).Value = "Não definida"
$ws.Cells.Item(11, 11).Value = ""
$ws.Cells.Item(11, 12).Value = "Não Binário"
$ws.Cells.Item(11, 13).NumberFormat = $dateFormat
$ws.Cells.Item(11, 13).Value = 37089
$ws.Cells.Item(11, 14).Value = "(19) 9 9877-1212"

# --- Row 12: antônio de mattos (fully processed) ---------------------------
$ws.Cells.Item(12, 1).Value = "antônio de mattos"
$ws.Cells.Item(12, 2).Value = ""
$ws.Cells.Item(12, 3).Value = "Sim"
$ws.Cells.Item(12, 4).Value = ""
$ws.Cells.Item(12, 5).Value = "Terça e Quinta"
$ws.Cells.Item(12, 6).Value = "16h00"
$ws.Cells.Item(12, 7).Value = "Jefferson"
$ws.Cells.Item(12, 8).Value = "Adulto B"
$ws.Cells.Item(12, 9).Value = 43
$ws.Cells.Item(12, 10).Value = "Não definida"
$ws.Cells.Item(12, 11).Value = ""
$ws.Cells.Item(12, 12).Value = "Masculino"
$ws.Cells.Item(12, 13).NumberFormat = $dateFormat
$ws.Cells.Item(12, 13).Value = 30294
$ws.Cells.Item(12, 14).Value = "(19) 9 9976-3211"

# --- Row 13: carla camuratti (added, not yet run through the age/category
#     logic - Idade stays blank and Categoria/Data keep their raw text) ----
$ws.Cells.Item(13, 1).Value = "carla camuratti"
$ws.Cells.Item(13, 2).Value = ""
$ws.Cells.Item(13, 3).Value = "Sim"
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 5).Value = "Terça e Quinta"
$ws.Cells.Item(13, 6).Value = "16h00"
$ws.Cells.Item(13, 7).Value = "Daniela"
$ws.Cells.Item(13, 8).Value = "Adulto A"
$ws.Cells.Item(13, 9).Value = ""
$ws.Cells.Item(13, 10).Value = "G50+"
$ws.Cells.Item(13, 11).Value = ""
$ws.Cells.Item(13, 12).Value = "Feminino"
# Force plain text (don't let COM auto-parse this as a date serial) and then
# drop the resulting quote-prefix style so the cell matches a bare inlineStr.
$ws.Cells.Item(13, 13).NumberFormat = "@"
$ws.Cells.Item(13, 13).Value = "05/09/1974"
$ws.Cells.Item(13, 13).ClearFormats()
$ws.Cells.Item(13, 14).Value = "(21) 9 9933-2876"
